$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010562896728516
$ws.Range("B1").Value = 2.123342990875244
$ws.Range("C1").Value = 6.112164497375488
$ws.Range("D1").Value = 1.324487924575806
$ws.Range("E1").Value = 1.290057063102722
